$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price and volume(1h) data per latest scrape
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.425.10'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.851.75'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.32'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.28%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4741'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2746'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.90%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06330'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +10.42%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.881.93'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07443'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.950'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '84.69'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6242'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.81%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.395.14'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.05%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '245.79'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +7.64%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.65'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007335'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.00%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.898'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.898'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '165.02'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.080'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '17.97'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.874'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1027'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.01%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -1.34%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.040'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.816'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04838'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.14%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.6971'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.710'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.87%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.01901'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +4.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.685'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.8779'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.63%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.992'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '106.55'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.18%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9999'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4064'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.501'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.157'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.49%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.27'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +6.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1196'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '33.93'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +4.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.581'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05503'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.349'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.29%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3694'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.01%  '
